$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3.0
$ws.Range("G2").Value = 8.723857666666666
$ws.Range("H2").Value = 26.171573
$ws.Range("I2").Value = 0.007098432040951201
$ws.Range("J2").Value = 0.007098432040951203
$ws.Range("K2").Value = 3.0
$ws.Range("M2").Value = 94.799851
$ws.Range("N2").Value = 284.399553
$ws.Range("O2").Value = 0.1681963571786457
$ws.Range("P2").Value = 0.1681963571786457
$ws.Range("Q2").Value = 827.0204069452076
$ws.Range("R2").Value = 7443.183662506869
$ws.Range("S2").Value = 0.001193930410968171
$ws.Range("T2").Value = 0.001193930410968171

# Row 3
$ws.Range("E3").Value = 3.0
$ws.Range("G3").Value = 8.723857666666666
$ws.Range("H3").Value = 26.171573
$ws.Range("I3").Value = 0.007098432040951201
$ws.Range("J3").Value = 0.007098432040951203
$ws.Range("K3").Value = 3.0
$ws.Range("M3").Value = 145.6413626666667
$ws.Range("N3").Value = 436.924088
$ws.Range("O3").Value = 0.2584006873076977
$ws.Range("P3").Value = 0.2584006873076978
$ws.Range("Q3").Value = 1270.55451828338
$ws.Range("R3").Value = 11434.99066455042
$ws.Range("S3").Value = 0.001834239718188774
$ws.Range("T3").Value = 0.001834239718188775

# Row 4
$ws.Range("E4").Value = 3.0
$ws.Range("G4").Value = 8.723857666666666
$ws.Range("H4").Value = 26.171573
$ws.Range("I4").Value = 0.007098432040951201
$ws.Range("J4").Value = 0.007098432040951203
$ws.Range("K4").Value = 3.0
$ws.Range("M4").Value = 72.77597066666667
$ws.Range("N4").Value = 218.327912
$ws.Range("O4").Value = 0.1291210168281099
$ws.Range("P4").Value = 0.1291210168281099
$ws.Range("Q4").Value = 634.8872096495083
$ws.Range("R4").Value = 5713.984886845576
$ws.Range("S4").Value = 0.0009165567630128542
$ws.Range("T4").Value = 0.0009165567630128545

# Row 5
$ws.Range("E5").Value = 3.0
$ws.Range("G5").Value = 8.723857666666666
$ws.Range("H5").Value = 26.171573
$ws.Range("I5").Value = 0.007098432040951201
$ws.Range("J5").Value = 0.007098432040951203
$ws.Range("K5").Value = 3.0
$ws.Range("M5").Value = 76.38610333333334
$ws.Range("N5").Value = 229.15831
$ws.Range("O5").Value = 0.1355262079445491
$ws.Range("P5").Value = 0.1355262079445491
$ws.Range("Q5").Value = 666.3814931912922
$ws.Range("R5").Value = 5997.43343872163
$ws.Range("S5").Value = 0.0009620235768622024
$ws.Range("T5").Value = 0.0009620235768622029

# Row 6
$ws.Range("E6").Value = 3.0
$ws.Range("G6").Value = 8.723857666666666
$ws.Range("H6").Value = 26.171573
$ws.Range("I6").Value = 0.007098432040951201
$ws.Range("J6").Value = 0.007098432040951203
$ws.Range("K6").Value = 3.0
$ws.Range("M6").Value = 58.18688599999999
$ws.Range("N6").Value = 174.560658
$ws.Range("O6").Value = 0.10323668399826
$ws.Range("P6").Value = 0.1032366839982601
$ws.Range("Q6").Value = 507.6141115305592
$ws.Range("R6").Value = 4568.527003775033
$ws.Range("S6").Value = 0.0007328185854948032
$ws.Range("T6").Value = 0.0007328185854948035

# Row 7
$ws.Range("E7").Value = 3.0
$ws.Range("G7").Value = 8.723857666666666
$ws.Range("H7").Value = 26.171573
$ws.Range("I7").Value = 0.007098432040951201
$ws.Range("J7").Value = 0.007098432040951203
$ws.Range("K7").Value = 3.0
$ws.Range("M7").Value = 115.8358916666667
$ws.Range("N7").Value = 347.5076749999999
$ws.Range("O7").Value = 0.2055190467427377
$ws.Range("P7").Value = 0.2055190467427377
$ws.Range("Q7").Value = 1010.535831591419
$ws.Range("R7").Value = 9094.822484322773
$ws.Range("S7").Value = 0.001458862986424397
$ws.Range("T7").Value = 0.001458862986424397

# Row 8
$ws.Range("E8").Value = 3.0
$ws.Range("G8").Value = 1065.000325333333
$ws.Range("H8").Value = 3195.000976
$ws.Range("I8").Value = 0.8665698962346957
$ws.Range("J8").Value = 0.8665698962346958
$ws.Range("K8").Value = 3.0
$ws.Range("M8").Value = 94.799851
$ws.Range("N8").Value = 284.399553
$ws.Range("O8").Value = 0.1681963571786457
$ws.Range("P8").Value = 0.1681963571786457
$ws.Range("Q8").Value = 100961.8721565515
$ws.Range("R8").Value = 908656.8494089639
$ws.Range("S8").Value = 0.1457538997873528
$ws.Range("T8").Value = 0.1457538997873528

# Row 9
$ws.Range("E9").Value = 3.0
$ws.Range("G9").Value = 1065.000325333333
$ws.Range("H9").Value = 3195.000976
$ws.Range("I9").Value = 0.8665698962346957
$ws.Range("J9").Value = 0.8665698962346958
$ws.Range("K9").Value = 3.0
$ws.Range("M9").Value = 145.6413626666667
$ws.Range("N9").Value = 436.924088
$ws.Range("O9").Value = 0.2584006873076977
$ws.Range("P9").Value = 0.2584006873076978
$ws.Range("Q9").Value = 155108.09862199
$ws.Range("R9").Value = 1395972.88759791
$ws.Range("S9").Value = 0.2239222567872057
$ws.Range("T9").Value = 0.2239222567872057

# Row 10
$ws.Range("E10").Value = 3.0
$ws.Range("G10").Value = 1065.000325333333
$ws.Range("H10").Value = 3195.000976
$ws.Range("I10").Value = 0.8665698962346957
$ws.Range("J10").Value = 0.8665698962346958
$ws.Range("K10").Value = 3.0
$ws.Range("M10").Value = 72.77597066666667
$ws.Range("N10").Value = 218.327912
$ws.Range("O10").Value = 0.1291210168281099
$ws.Range("P10").Value = 0.1291210168281099
$ws.Range("Q10").Value = 77506.43243644913
$ws.Range("R10").Value = 697557.8919280422
$ws.Range("S10").Value = 0.1118923861544536
$ws.Range("T10").Value = 0.1118923861544536

# Row 11
$ws.Range("E11").Value = 3.0
$ws.Range("G11").Value = 1065.000325333333
$ws.Range("H11").Value = 3195.000976
$ws.Range("I11").Value = 0.8665698962346957
$ws.Range("J11").Value = 0.8665698962346958
$ws.Range("K11").Value = 3.0
$ws.Range("M11").Value = 76.38610333333334
$ws.Range("N11").Value = 229.15831
$ws.Range("O11").Value = 0.1355262079445491
$ws.Range("P11").Value = 0.1355262079445491
$ws.Range("Q11").Value = 81351.22490094563
$ws.Range("R11").Value = 732161.0241085107
$ws.Range("S11").Value = 0.1174429319555897
$ws.Range("T11").Value = 0.1174429319555897

# Row 12
$ws.Range("E12").Value = 3.0
$ws.Range("G12").Value = 1065.000325333333
$ws.Range("H12").Value = 3195.000976
$ws.Range("I12").Value = 0.8665698962346957
$ws.Range("J12").Value = 0.8665698962346958
$ws.Range("K12").Value = 3.0
$ws.Range("M12").Value = 58.18688599999999
$ws.Range("N12").Value = 174.560658
$ws.Range("O12").Value = 0.10323668399826
$ws.Range("P12").Value = 0.1032366839982601
$ws.Range("Q12").Value = 61969.05252013358
$ws.Range("R12").Value = 557721.4726812022
$ws.Range("S12").Value = 0.08946180253998627
$ws.Range("T12").Value = 0.0894618025399863

# Row 13
$ws.Range("E13").Value = 3.0
$ws.Range("G13").Value = 1065.000325333333
$ws.Range("H13").Value = 3195.000976
$ws.Range("I13").Value = 0.8665698962346957
$ws.Range("J13").Value = 0.8665698962346958
$ws.Range("K13").Value = 3.0
$ws.Range("M13").Value = 115.8358916666667
$ws.Range("N13").Value = 347.5076749999999
$ws.Range("O13").Value = 0.2055190467427377
$ws.Range("P13").Value = 0.2055190467427377
$ws.Range("Q13").Value = 123365.2623102768
$ws.Range("R13").Value = 1110287.360792491
$ws.Range("S13").Value = 0.1780966190101078
$ws.Range("T13").Value = 0.1780966190101078

# Row 14
$ws.Range("E14").Value = 3.0
$ws.Range("G14").Value = 0.8450703333333333
$ws.Range("H14").Value = 2.535211
$ws.Range("I14").Value = 0.0006876171712327699
$ws.Range("J14").Value = 0.0006876171712327699
$ws.Range("K14").Value = 3.0
$ws.Range("M14").Value = 94.799851
$ws.Range("N14").Value = 284.399553
$ws.Range("O14").Value = 0.1681963571786457
$ws.Range("P14").Value = 0.1681963571786457
$ws.Range("Q14").Value = 80.11254168452032
$ws.Range("R14").Value = 721.012875160683
$ws.Range("S14").Value = 0.0001156547033348369
$ws.Range("T14").Value = 0.0001156547033348369

# Row 15
$ws.Range("E15").Value = 3.0
$ws.Range("G15").Value = 0.8450703333333333
$ws.Range("H15").Value = 2.535211
$ws.Range("I15").Value = 0.0006876171712327699
$ws.Range("J15").Value = 0.0006876171712327699
$ws.Range("K15").Value = 3.0
$ws.Range("M15").Value = 145.6413626666667
$ws.Range("N15").Value = 436.924088
$ws.Range("O15").Value = 0.2584006873076977
$ws.Range("P15").Value = 0.2584006873076978
$ws.Range("Q15").Value = 123.0771948958409
$ws.Range("R15").Value = 1107.694754062568
$ws.Range("S15").Value = 0.0001776807496511226
$ws.Range("T15").Value = 0.0001776807496511226

# Row 16
$ws.Range("E16").Value = 3.0
$ws.Range("G16").Value = 0.8450703333333333
$ws.Range("H16").Value = 2.535211
$ws.Range("I16").Value = 0.0006876171712327699
$ws.Range("J16").Value = 0.0006876171712327699
$ws.Range("K16").Value = 3.0
$ws.Range("M16").Value = 72.77597066666667
$ws.Range("N16").Value = 218.327912
$ws.Range("O16").Value = 0.1291210168281099
$ws.Range("P16").Value = 0.1291210168281099
$ws.Range("Q16").Value = 61.50081378993688
$ws.Range("R16").Value = 553.5073241094319
$ws.Range("S16").Value = 0.00008878582833804378
$ws.Range("T16").Value = 0.00008878582833804378

# Row 17
$ws.Range("E17").Value = 3.0
$ws.Range("G17").Value = 0.8450703333333333
$ws.Range("H17").Value = 2.535211
$ws.Range("I17").Value = 0.0006876171712327699
$ws.Range("J17").Value = 0.0006876171712327699
$ws.Range("K17").Value = 3.0
$ws.Range("M17").Value = 76.38610333333334
$ws.Range("N17").Value = 229.15831
$ws.Range("O17").Value = 0.1355262079445491
$ws.Range("P17").Value = 0.1355262079445491
$ws.Range("Q17").Value = 64.55162980593444
$ws.Range("R17").Value = 580.9646682534101
$ws.Range("S17").Value = 0.00009319014773473499
$ws.Range("T17").Value = 0.000093190147734735

# Row 18
$ws.Range("E18").Value = 3.0
$ws.Range("G18").Value = 0.8450703333333333
$ws.Range("H18").Value = 2.535211
$ws.Range("I18").Value = 0.0006876171712327699
$ws.Range("J18").Value = 0.0006876171712327699
$ws.Range("K18").Value = 3.0
$ws.Range("M18").Value = 58.18688599999999
$ws.Range("N18").Value = 174.560658
$ws.Range("O18").Value = 0.10323668399826
$ws.Range("P18").Value = 0.1032366839982601
$ws.Range("Q18").Value = 49.17201114764866
$ws.Range("R18").Value = 442.548100328838
$ws.Range("S18").Value = 0.00007098731661833493
$ws.Range("T18").Value = 0.00007098731661833493

# Row 19
$ws.Range("E19").Value = 3.0
$ws.Range("G19").Value = 0.8450703333333333
$ws.Range("H19").Value = 2.535211
$ws.Range("I19").Value = 0.0006876171712327699
$ws.Range("J19").Value = 0.0006876171712327699
$ws.Range("K19").Value = 3.0
$ws.Range("M19").Value = 115.8358916666667
$ws.Range("N19").Value = 347.5076749999999
$ws.Range("O19").Value = 0.2055190467427377
$ws.Range("P19").Value = 0.2055190467427377
$ws.Range("Q19").Value = 97.88947558271387
$ws.Range("R19").Value = 881.0052802444249
$ws.Range("S19").Value = 0.0001413184255556967
$ws.Range("T19").Value = 0.0001413184255556967

# Row 20
$ws.Range("E20").Value = 3.0
$ws.Range("G20").Value = 1.780731666666667
$ws.Range("H20").Value = 5.342195
$ws.Range("I20").Value = 0.001448946464051256
$ws.Range("J20").Value = 0.001448946464051256
$ws.Range("K20").Value = 3.0
$ws.Range("M20").Value = 94.799851
$ws.Range("N20").Value = 284.399553
$ws.Range("O20").Value = 0.1681963571786457
$ws.Range("P20").Value = 0.1681963571786457
$ws.Range("Q20").Value = 168.8130966709817
$ws.Range("R20").Value = 1519.317870038835
$ws.Range("S20").Value = 0.0002437075170003007
$ws.Range("T20").Value = 0.0002437075170003007

# Row 21
$ws.Range("E21").Value = 3.0
$ws.Range("G21").Value = 1.780731666666667
$ws.Range("H21").Value = 5.342195
$ws.Range("I21").Value = 0.001448946464051256
$ws.Range("J21").Value = 0.001448946464051256
$ws.Range("K21").Value = 3.0
$ws.Range("M21").Value = 145.6413626666667
$ws.Range("N21").Value = 436.924088
$ws.Range("O21").Value = 0.2584006873076977
$ws.Range("P21").Value = 0.2584006873076978
$ws.Range("Q21").Value = 259.3481864770177
$ws.Range("R21").Value = 2334.13367829316
$ws.Range("S21").Value = 0.0003744087621829027
$ws.Range("T21").Value = 0.0003744087621829028

# Row 22
$ws.Range("E22").Value = 3.0
$ws.Range("G22").Value = 1.780731666666667
$ws.Range("H22").Value = 5.342195
$ws.Range("I22").Value = 0.001448946464051256
$ws.Range("J22").Value = 0.001448946464051256
$ws.Range("K22").Value = 3.0
$ws.Range("M22").Value = 72.77597066666667
$ws.Range("N22").Value = 218.327912
$ws.Range("O22").Value = 0.1291210168281099
$ws.Range("P22").Value = 0.1291210168281099
$ws.Range("Q22").Value = 129.5944755385378
$ws.Range("R22").Value = 1166.35027984684
$ws.Range("S22").Value = 0.0001870894407677924
$ws.Range("T22").Value = 0.0001870894407677924

# Row 23
$ws.Range("E23").Value = 3.0
$ws.Range("G23").Value = 1.780731666666667
$ws.Range("H23").Value = 5.342195
$ws.Range("I23").Value = 0.001448946464051256
$ws.Range("J23").Value = 0.001448946464051256
$ws.Range("K23").Value = 3.0
$ws.Range("M23").Value = 76.38610333333334
$ws.Range("N23").Value = 229.15831
$ws.Range("O23").Value = 0.1355262079445491
$ws.Range("P23").Value = 0.1355262079445491
$ws.Range("Q23").Value = 136.0231530989389
$ws.Range("R23").Value = 1224.20837789045
$ws.Range("S23").Value = 0.0001963702197875296
$ws.Range("T23").Value = 0.0001963702197875296

# Row 24
$ws.Range("E24").Value = 3.0
$ws.Range("G24").Value = 1.780731666666667
$ws.Range("H24").Value = 5.342195
$ws.Range("I24").Value = 0.001448946464051256
$ws.Range("J24").Value = 0.001448946464051256
$ws.Range("K24").Value = 3.0
$ws.Range("M24").Value = 58.18688599999999
$ws.Range("N24").Value = 174.560658
$ws.Range("O24").Value = 0.10323668399826
$ws.Range("P24").Value = 0.1032366839982601
$ws.Range("Q24").Value = 103.6152304849233
$ws.Range("R24").Value = 932.53707436431
$ws.Range("S24").Value = 0.0001495844282396557
$ws.Range("T24").Value = 0.0001495844282396557

# Row 25
$ws.Range("E25").Value = 3.0
$ws.Range("G25").Value = 1.780731666666667
$ws.Range("H25").Value = 5.342195
$ws.Range("I25").Value = 0.001448946464051256
$ws.Range("J25").Value = 0.001448946464051256
$ws.Range("K25").Value = 3.0
$ws.Range("M25").Value = 115.8358916666667
$ws.Range("N25").Value = 347.5076749999999
$ws.Range("O25").Value = 0.2055190467427377
$ws.Range("P25").Value = 0.2055190467427377
$ws.Range("Q25").Value = 206.2726404274028
$ws.Range("R25").Value = 1856.453763846625
$ws.Range("S25").Value = 0.0002977860960730745
$ws.Range("T25").Value = 0.0002977860960730745

# Row 26
$ws.Range("E26").Value = 3.0
$ws.Range("G26").Value = 0.2508443333333333
$ws.Range("H26").Value = 0.752533
$ws.Range("I26").Value = 0.0002041071187839237
$ws.Range("J26").Value = 0.0002041071187839237
$ws.Range("K26").Value = 3.0
$ws.Range("M26").Value = 94.799851
$ws.Range("N26").Value = 284.399553
$ws.Range("O26").Value = 0.1681963571786457
$ws.Range("P26").Value = 0.1681963571786457
$ws.Range("Q26").Value = 23.78000542419434
$ws.Range("R26").Value = 214.020048817749
$ws.Range("S26").Value = 0.00003433007385368509
$ws.Range("T26").Value = 0.0000343300738536851

# Row 27
$ws.Range("E27").Value = 3.0
$ws.Range("G27").Value = 0.2508443333333333
$ws.Range("H27").Value = 0.752533
$ws.Range("I27").Value = 0.0002041071187839237
$ws.Range("J27").Value = 0.0002041071187839237
$ws.Range("K27").Value = 3.0
$ws.Range("M27").Value = 145.6413626666667
$ws.Range("N27").Value = 436.924088
$ws.Range("O27").Value = 0.2584006873076977
$ws.Range("P27").Value = 0.2584006873076978
$ws.Range("Q27").Value = 36.53331052387822
$ws.Range("R27").Value = 328.799794714904
$ws.Range("S27").Value = 0.00005274141977815979
$ws.Range("T27").Value = 0.0000527414197781598

# Row 28
$ws.Range("E28").Value = 3.0
$ws.Range("G28").Value = 0.2508443333333333
$ws.Range("H28").Value = 0.752533
$ws.Range("I28").Value = 0.0002041071187839237
$ws.Range("J28").Value = 0.0002041071187839237
$ws.Range("K28").Value = 3.0
$ws.Range("M28").Value = 72.77597066666667
$ws.Range("N28").Value = 218.327912
$ws.Range("O28").Value = 0.1291210168281099
$ws.Range("P28").Value = 0.1291210168281099
$ws.Range("Q28").Value = 18.25543984456622
$ws.Range("R28").Value = 164.298958601096
$ws.Range("S28").Value = 0.00002635451871923603
$ws.Range("T28").Value = 0.00002635451871923603

# Row 29
$ws.Range("E29").Value = 3.0
$ws.Range("G29").Value = 0.2508443333333333
$ws.Range("H29").Value = 0.752533
$ws.Range("I29").Value = 0.0002041071187839237
$ws.Range("J29").Value = 0.0002041071187839237
$ws.Range("K29").Value = 3.0
$ws.Range("M29").Value = 76.38610333333334
$ws.Range("N29").Value = 229.15831
$ws.Range("O29").Value = 0.1355262079445491
$ws.Range("P29").Value = 0.1355262079445491
$ws.Range("Q29").Value = 19.16102116658111
$ws.Range("R29").Value = 172.44919049923
$ws.Range("S29").Value = 0.00002766186382327283
$ws.Range("T29").Value = 0.00002766186382327284

# Row 30
$ws.Range("E30").Value = 3.0
$ws.Range("G30").Value = 0.2508443333333333
$ws.Range("H30").Value = 0.752533
$ws.Range("I30").Value = 0.0002041071187839237
$ws.Range("J30").Value = 0.0002041071187839237
$ws.Range("K30").Value = 3.0
$ws.Range("M30").Value = 58.18688599999999
$ws.Range("N30").Value = 174.560658
$ws.Range("O30").Value = 0.10323668399826
$ws.Range("P30").Value = 0.1032366839982601
$ws.Range("Q30").Value = 14.59585062741266
$ws.Range("R30").Value = 131.362655646714
$ws.Range("S30").Value = 0.00002107134212369126
$ws.Range("T30").Value = 0.00002107134212369126

# Row 31
$ws.Range("E31").Value = 3.0
$ws.Range("G31").Value = 0.2508443333333333
$ws.Range("H31").Value = 0.752533
$ws.Range("I31").Value = 0.0002041071187839237
$ws.Range("J31").Value = 0.0002041071187839237
$ws.Range("K31").Value = 3.0
$ws.Range("M31").Value = 115.8358916666667
$ws.Range("N31").Value = 347.5076749999999
$ws.Range("O31").Value = 0.2055190467427377
$ws.Range("P31").Value = 0.2055190467427377
$ws.Range("Q31").Value = 29.05677702119722
$ws.Range("R31").Value = 261.510993190775
$ws.Range("S31").Value = 0.00004194790048587874
$ws.Range("T31").Value = 0.00004194790048587874

# Row 32
$ws.Range("E32").Value = 3.0
$ws.Range("G32").Value = 152.382926
$ws.Range("H32").Value = 457.148778
$ws.Range("I32").Value = 0.1239910009702851
$ws.Range("J32").Value = 0.1239910009702851
$ws.Range("K32").Value = 3.0
$ws.Range("M32").Value = 94.799851
$ws.Range("N32").Value = 284.399553
$ws.Range("O32").Value = 0.1681963571786457
$ws.Range("P32").Value = 0.1681963571786457
$ws.Range("Q32").Value = 14445.87867974403
$ws.Range("R32").Value = 130012.9081176962
$ws.Range("S32").Value = 0.02085483468613588
$ws.Range("T32").Value = 0.02085483468613588

# Row 33
$ws.Range("E33").Value = 3.0
$ws.Range("G33").Value = 152.382926
$ws.Range("H33").Value = 457.148778
$ws.Range("I33").Value = 0.1239910009702851
$ws.Range("J33").Value = 0.1239910009702851
$ws.Range("K33").Value = 3.0
$ws.Range("M33").Value = 145.6413626666667
$ws.Range("N33").Value = 436.924088
$ws.Range("O33").Value = 0.2584006873076977
$ws.Range("P33").Value = 0.2584006873076978
$ws.Range("Q33").Value = 22193.25698977382
$ws.Range("R33").Value = 199739.3129079645
$ws.Range("S33").Value = 0.03203935987069109
$ws.Range("T33").Value = 0.0320393598706911

# Row 34
$ws.Range("E34").Value = 3.0
$ws.Range("G34").Value = 152.382926
$ws.Range("H34").Value = 457.148778
$ws.Range("I34").Value = 0.1239910009702851
$ws.Range("J34").Value = 0.1239910009702851
$ws.Range("K34").Value = 3.0
$ws.Range("M34").Value = 72.77597066666667
$ws.Range("N34").Value = 218.327912
$ws.Range("O34").Value = 0.1291210168281099
$ws.Range("P34").Value = 0.1291210168281099
$ws.Range("Q34").Value = 11089.81535267684
$ws.Range("R34").Value = 99808.33817409154
$ws.Range("S34").Value = 0.01600984412281837
$ws.Range("T34").Value = 0.01600984412281837

# Row 35
$ws.Range("E35").Value = 3.0
$ws.Range("G35").Value = 152.382926
$ws.Range("H35").Value = 457.148778
$ws.Range("I35").Value = 0.1239910009702851
$ws.Range("J35").Value = 0.1239910009702851
$ws.Range("K35").Value = 3.0
$ws.Range("M35").Value = 76.38610333333334
$ws.Range("N35").Value = 229.15831
$ws.Range("O35").Value = 0.1355262079445491
$ws.Range("P35").Value = 0.1355262079445491
$ws.Range("Q35").Value = 11639.93793167169
$ws.Range("R35").Value = 104759.4413850452
$ws.Range("S35").Value = 0.01680403018075165
$ws.Range("T35").Value = 0.01680403018075165

# Row 36
$ws.Range("E36").Value = 3.0
$ws.Range("G36").Value = 152.382926
$ws.Range("H36").Value = 457.148778
$ws.Range("I36").Value = 0.1239910009702851
$ws.Range("J36").Value = 0.1239910009702851
$ws.Range("K36").Value = 3.0
$ws.Range("M36").Value = 58.18688599999999
$ws.Range("N36").Value = 174.560658
$ws.Range("O36").Value = 0.10323668399826
$ws.Range("P36").Value = 0.1032366839982601
$ws.Range("Q36").Value = 8866.687943508436
$ws.Range("R36").Value = 79800.19149157591
$ws.Range("S36").Value = 0.01280041978579728
$ws.Range("T36").Value = 0.01280041978579728

# Row 37
$ws.Range("E37").Value = 3.0
$ws.Range("G37").Value = 152.382926
$ws.Range("H37").Value = 457.148778
$ws.Range("I37").Value = 0.1239910009702851
$ws.Range("J37").Value = 0.1239910009702851
$ws.Range("K37").Value = 3.0
$ws.Range("M37").Value = 115.8358916666667
$ws.Range("N37").Value = 347.5076749999999
$ws.Range("O37").Value = 0.2055190467427377
$ws.Range("P37").Value = 0.2055190467427377
$ws.Range("Q37").Value = 17651.41210798568
$ws.Range("R37").Value = 158862.7089718711
$ws.Range("S37").Value = 0.02548251232409086
$ws.Range("T37").Value = 0.02548251232409087
